$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Narea_model")

# 1. Update header text from "RelImp" to "Relative Importance"
$ws.Range("E1").Value = "Relative Importance"

# 2. Give the Relative-Importance data column (E2:E16) its own number
#    format (0.00), which creates a new cell style distinct from the
#    other "0.000"-formatted numeric columns.
$ws.Range("E2:E16").NumberFormat = "0.00"

# 3. Widen column E so the new, longer header text fits.
$ws.Columns("E:E").ColumnWidth = 17

# 4. Update the active selection left on the sheet.
$ws.Range("E19").Select()
